# Scheduled refresh of live market-price driven columns (H:N) across all
# Leve-profit tracking sheets. Values below are the latest snapshot pulled
# from the market data source for each Leve row (looked up by sheet + row).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2: Mercury Rising
$ws.Range("H2").Value = 176.66667
$ws.Range("I2").Value = 176.66667
$ws.Range("K2").Value = 176.66667
$ws.Range("M2").Value = -63.66667000000001

# Row 18: You Grow, Girl
$ws.Range("H18").Value = 1865.1666
$ws.Range("I18").Value = 1737.8
$ws.Range("J18").Value = 2502
$ws.Range("K18").Value = 1737.8
$ws.Range("L18").Value = 2502
$ws.Range("M18").Value = -1453.8
$ws.Range("N18").Value = -3070

# Row 48: The Sting of Conscience
$ws.Range("H48").Value = 3000
$ws.Range("J48").Value = 3000
$ws.Range("L48").Value = 9000
$ws.Range("N48").Value = -9584

# Row 56: Sleepless in Silvertear
$ws.Range("H56").Value = 3000
$ws.Range("J56").Value = 3000
$ws.Range("L56").Value = 9000
$ws.Range("N56").Value = -10068

# Row 76: Warding Off Temptation
$ws.Range("H76").Value = 5244.2593
$ws.Range("J76").Value = 6268.7334
$ws.Range("L76").Value = 6268.7334
$ws.Range("N76").Value = -6898.7334

# Row 79: The Garden of Arcane Delights (L)
$ws.Range("H79").Value = 5244.2593
$ws.Range("J79").Value = 6268.7334
$ws.Range("L79").Value = 6268.7334
$ws.Range("N79").Value = -8452.733400000001

# Row 107: Another Man's Ink
$ws.Range("H107").Value = 2336.879
$ws.Range("I107").Value = 2647.9644
$ws.Range("J107").Value = 594.8
$ws.Range("K107").Value = 2647.9644
$ws.Range("L107").Value = 594.8
$ws.Range("M107").Value = -727.9643999999998
$ws.Range("N107").Value = -4434.8

# Row 132: Fast-forwarding Flora
$ws.Range("H132").Value = 20844
$ws.Range("I132").Value = 29264.857
$ws.Range("K132").Value = 87794.571
$ws.Range("M132").Value = -85264.571

# Row 138: All-night Crafting
$ws.Range("H138").Value = 4291.788
$ws.Range("I138").Value = 3861.5
$ws.Range("J138").Value = 5890
$ws.Range("K138").Value = 11584.5
$ws.Range("L138").Value = 17670
$ws.Range("M138").Value = -6444.5
$ws.Range("N138").Value = -27950

# Row 141: Remedy for Reason
$ws.Range("H141").Value = 4746.0303
$ws.Range("I141").Value = 2594
$ws.Range("J141").Value = 11471.125
$ws.Range("K141").Value = 7782
$ws.Range("L141").Value = 34413.375
$ws.Range("M141").Value = -2602
$ws.Range("N141").Value = -44773.375

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 238835.86
$ws.Range("I32").Value = 242141.47
$ws.Range("K32").Value = 242141.47
$ws.Range("M32").Value = -241854.47

# Row 45: Hollow Hallmarks
$ws.Range("H45").Value = 3424.875
$ws.Range("J45").Value = 5133.3335
$ws.Range("L45").Value = 5133.3335
$ws.Range("N45").Value = -5887.3335

# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 3983.5715
$ws.Range("I61").Value = 3314.1667
$ws.Range("K61").Value = 3314.1667
$ws.Range("M61").Value = -3102.1667

# Row 74: As the Bolt Flies
$ws.Range("H74").Value = 6956.1914
$ws.Range("I74").Value = 3837.0527
$ws.Range("J74").Value = 20125.889
$ws.Range("K74").Value = 3837.0527
$ws.Range("L74").Value = 20125.889
$ws.Range("M74").Value = -2963.0527
$ws.Range("N74").Value = -21873.889

# Row 77: Heavy Metal Banned (L)
$ws.Range("H77").Value = 6956.1914
$ws.Range("I77").Value = 3837.0527
$ws.Range("J77").Value = 20125.889
$ws.Range("K77").Value = 19185.2635
$ws.Range("L77").Value = 100629.445
$ws.Range("M77").Value = -14817.2635
$ws.Range("N77").Value = -109365.445

# Row 136: Metal with Mettle
$ws.Range("H136").Value = 3983.5715
$ws.Range("I136").Value = 3314.1667
$ws.Range("K136").Value = 9942.500100000001
$ws.Range("M136").Value = -7392.500100000001

$ws = $wb.Worksheets.Item("BSM")
# Row 22: Riveting Run
$ws.Range("H22").Value = 312.5
$ws.Range("I22").Value = 250
$ws.Range("K22").Value = 250
$ws.Range("M22").Value = -77

# Row 86: Through Thick and Thin
$ws.Range("H86").Value = 2596
$ws.Range("I86").Value = 2399.3333
$ws.Range("K86").Value = 2399.3333
$ws.Range("M86").Value = -1276.3333

# Row 89: Piercing Eyes Deserve Piercing Shafts (L)
$ws.Range("H89").Value = 2596
$ws.Range("I89").Value = 2399.3333
$ws.Range("K89").Value = 11996.6665
$ws.Range("M89").Value = -6380.666499999999

# Row 105: Ingot to Wing It
$ws.Range("H105").Value = 6002.1113
$ws.Range("I105").Value = 2007.1428
$ws.Range("K105").Value = 2007.1428
$ws.Range("M105").Value = -260.1428000000001

# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 7166.3335
$ws.Range("I134").Value = 7399.6
$ws.Range("K134").Value = 22198.8
$ws.Range("M134").Value = -19663.8

$ws = $wb.Worksheets.Item("CRP")
# Row 86: Birch, Please
$ws.Range("H86").Value = 95775.3
$ws.Range("J86").Value = 14444
$ws.Range("L86").Value = 14444
$ws.Range("N86").Value = -16690

# Row 87: Anatomy of a Drill Bit
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

# Row 89: Built This City on Blocks and Soul (L)
$ws.Range("H89").Value = 95775.3
$ws.Range("J89").Value = 14444
$ws.Range("L89").Value = 72220
$ws.Range("N89").Value = -83452

# Row 90: Pulling Them to the Grind (L)
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

# Row 107: Built to Last
$ws.Range("H107").Value = 1068.8
$ws.Range("I107").Value = 790
$ws.Range("K107").Value = 790
$ws.Range("M107").Value = 1130

# Row 122: Timber of Tenkonto
$ws.Range("H122").Value = 13320.565
$ws.Range("I122").Value = 2387.95
$ws.Range("K122").Value = 7163.849999999999
$ws.Range("M122").Value = -4713.849999999999

# Row 132: Hull Lotta Damage
$ws.Range("H132").Value = 2145.2856
$ws.Range("I132").Value = 2148.8293
$ws.Range("K132").Value = 6446.4879
$ws.Range("M132").Value = -3916.4879

# Row 134: Wood You Be Quiet
$ws.Range("H134").Value = 3211.353
$ws.Range("I134").Value = 2964
$ws.Range("J134").Value = 3664.8333
$ws.Range("K134").Value = 8892
$ws.Range("L134").Value = 10994.4999
$ws.Range("M134").Value = -6357
$ws.Range("N134").Value = -16064.4999

$ws = $wb.Worksheets.Item("CUL")
# Row 9: Jack of All Plates
$ws.Range("H9").Value = 558672.5
$ws.Range("J9").Value = 542910.8
$ws.Range("L9").Value = 1628732.4
$ws.Range("N9").Value = -1629180.4

# Row 40: True Grits
$ws.Range("H40").Value = 56.5
$ws.Range("I40").Value = 50
$ws.Range("J40").Value = 102
$ws.Range("K40").Value = 200
$ws.Range("L40").Value = 408
$ws.Range("M40").Value = -131
$ws.Range("N40").Value = -546

# Row 131: The Mountain Steeped
$ws.Range("H131").Value = 2357.875
$ws.Range("I131").Value = 1074.75
$ws.Range("J131").Value = 2474.5227
$ws.Range("K131").Value = 3224.25
$ws.Range("L131").Value = 7423.5681
$ws.Range("M131").Value = 1815.75
$ws.Range("N131").Value = -17503.5681

$ws = $wb.Worksheets.Item("GSM")
# Row 97: If I'd a Koppranickel for Every Time...
$ws.Range("H97").Value = 1082.1177
$ws.Range("I97").Value = 1238.1538
$ws.Range("K97").Value = 1238.1538
$ws.Range("M97").Value = -742.1538

# Row 102: Put the Metal to the Peddle
$ws.Range("H102").Value = 3042.111
$ws.Range("I102").Value = 3042.111
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 3042.111
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -1420.111
$ws.Range("N102").ClearContents()

# Row 113: Copious Crystal Cannons
$ws.Range("H113").Value = 2178.6843

# Row 132: On Board for Lar
$ws.Range("H132").Value = 21779.422
$ws.Range("I132").Value = 22930.883
$ws.Range("K132").Value = 68792.649
$ws.Range("M132").Value = -66262.649

$ws = $wb.Worksheets.Item("LTW")
# Row 22: Skin off Their Backs
$ws.Range("H22").Value = 3166.0557
$ws.Range("I22").Value = 2799.4
$ws.Range("J22").Value = 3307.077
$ws.Range("K22").Value = 2799.4
$ws.Range("L22").Value = 3307.077
$ws.Range("M22").Value = -2504.4
$ws.Range("N22").Value = -3897.077

# Row 27: Fire and Hide
$ws.Range("H27").Value = 3166.0557
$ws.Range("I27").Value = 2799.4
$ws.Range("J27").Value = 3307.077
$ws.Range("K27").Value = 2799.4
$ws.Range("L27").Value = 3307.077
$ws.Range("M27").Value = -2692.4
$ws.Range("N27").Value = -3521.077

# Row 46: Supply Side Logic
$ws.Range("H46").Value = 3764.5806
$ws.Range("I46").Value = 2471.9
$ws.Range("J46").Value = 4380.143
$ws.Range("K46").Value = 2471.9
$ws.Range("L46").Value = 4380.143
$ws.Range("M46").Value = -2283.9
$ws.Range("N46").Value = -4756.143

# Row 68: You Could Say It's a Moving Target
$ws.Range("H68").Value = 10018.625
$ws.Range("I68").Value = 8580.111000000001
$ws.Range("J68").Value = 14334.167
$ws.Range("K68").Value = 8580.111000000001
$ws.Range("L68").Value = 14334.167
$ws.Range("M68").Value = -7831.111000000001
$ws.Range("N68").Value = -15832.167

# Row 71: They Call It Bloody Mary (L)
$ws.Range("H71").Value = 10018.625
$ws.Range("I71").Value = 8580.111000000001
$ws.Range("J71").Value = 14334.167
$ws.Range("K71").Value = 42900.55500000001
$ws.Range("L71").Value = 71670.83499999999
$ws.Range("M71").Value = -39156.55500000001
$ws.Range("N71").Value = -79158.83499999999

# Row 93: Hide to Go Seek
$ws.Range("H93").Value = 2128.375
$ws.Range("I93").Value = 1081.9286
$ws.Range("K93").Value = 1081.9286
$ws.Range("M93").Value = 166.0714

# Row 132: Tenets of Tanning
$ws.Range("H132").Value = 2694.318
$ws.Range("I132").Value = 2011
$ws.Range("J132").Value = 3890.125
$ws.Range("K132").Value = 6033
$ws.Range("L132").Value = 11670.375
$ws.Range("M132").Value = -3503
$ws.Range("N132").Value = -16730.375

$ws = $wb.Worksheets.Item("WVR")
# Row 126: A Polished Purchase
$ws.Range("H126").Value = 2150.9
$ws.Range("I126").Value = 1700.5714
$ws.Range("K126").Value = 5101.7142
$ws.Range("M126").Value = -2631.7142

# Row 132: Comfy Cabins
$ws.Range("H132").Value = 2751.0386
$ws.Range("I132").Value = 1974.4615
$ws.Range("J132").Value = 3527.6155
$ws.Range("K132").Value = 5923.3845
$ws.Range("L132").Value = 10582.8465
$ws.Range("M132").Value = -3393.3845
$ws.Range("N132").Value = -15642.8465
